$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Contenu du stage" block (rows 16-23): update real counts + real
# percentages (previously every row but C# showed placeholder 0 %/100 %).

# Helper: write a cell as plain literal TEXT (shared string, no style /
# number-format changes), mirroring how the original "xx %" cells are
# stored (t="s", no cell style). A direct .Value = "73.47 %" assignment
# gets auto-detected by Excel as a percentage number, so instead we set a
# text-returning formula and immediately paste-special just the value
# back over itself, which collapses it to a literal string cell.
function Set-TextValue {
    param($addr, [string]$text)
    $r = $ws.Range($addr)
    $escaped = $text.Replace('"', '""')
    $r.Formula = '="' + $escaped + '"'
    $r.Copy()
    $r.PasteSpecial(-4163)
}

# Counts (column E)
$ws.Range("E16").Value = 8
$ws.Range("E17").Value = 36
$ws.Range("E19").Value = 2
$ws.Range("E20").Value = 3

# Percentages (column G)
Set-TextValue "G16" "16.33 %"
Set-TextValue "G17" "73.47 %"
Set-TextValue "G19" "4.08 %"
Set-TextValue "G20" "6.12 %"

$wb.Save()
